# "fixes for WO and CPQ"
# The "Routing Master" sheet keeps a sample/reference row (row 2) that
# mirrors a Salesforce Engineering Item ("Item Number" + its record "Id").
# That sample pair is stale; point it at the current Engineering Item
# record instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

# B2: Item Number, D2: Id - update the sample row to the new Pro-PEItem record
$ws.Range("B2").Value = "Pro-PEItem-F951G"
$ws.Range("D2").Value = "a345f000000uUJmAAM"

# The new values are longer than the old ones, so re-fit the two columns
# (both columns were already best-fit / auto-sized).
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()
